$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 4731.5586
$ws.Range("I15").Value = 4731.5586
$ws.Range("K15").Value = 14194.6758
$ws.Range("M15").Value = -14025.6758
$ws.Range("H32").Value = 2033.3334
$ws.Range("J32").Value = 2375
$ws.Range("L32").Value = 2375
$ws.Range("N32").Value = -3027
$ws.Range("H51").Value = 4166.6665
$ws.Range("J51").Value = 4166.6665
$ws.Range("L51").Value = 4166.6665
$ws.Range("N51").Value = -5134.6665
$ws.Range("H116").Value = 3364.1428
$ws.Range("I116").Value = 3274.875
$ws.Range("K116").Value = 3274.875
$ws.Range("M116").Value = 167.125
$ws.Range("H127").Value = 1110.5555
$ws.Range("J127").Value = 1184
$ws.Range("L127").Value = 3552
$ws.Range("N127").Value = -13472
$ws.Range("H137").Value = 1287.0625
$ws.Range("I137").Value = 1276.1538
$ws.Range("J137").Value = 1334.3334
$ws.Range("K137").Value = 3828.4614
$ws.Range("L137").Value = 4003.0002
$ws.Range("M137").Value = -1278.4614
$ws.Range("N137").Value = -9103.0002
$ws.Range("H138").Value = 1374.3062
$ws.Range("I138").Value = 694.0625
$ws.Range("J138").Value = 1704.1212
$ws.Range("K138").Value = 2082.1875
$ws.Range("L138").Value = 5112.363600000001
$ws.Range("M138").Value = 3057.8125
$ws.Range("N138").Value = -15392.3636

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2173.5
$ws.Range("I61").Value = 1841.875
$ws.Range("K61").Value = 1841.875
$ws.Range("M61").Value = -1629.875
$ws.Range("H74").Value = 2519.5715
$ws.Range("I74").Value = 1027.25
$ws.Range("K74").Value = 1027.25
$ws.Range("M74").Value = -153.25
$ws.Range("H77").Value = 2519.5715
$ws.Range("I77").Value = 1027.25
$ws.Range("K77").Value = 5136.25
$ws.Range("M77").Value = -768.25
$ws.Range("H102").Value = 8334205
$ws.Range("I102").Value = 9260151
$ws.Range("K102").Value = 9260151
$ws.Range("M102").Value = -9258529
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
$ws.Range("H136").Value = 2173.5
$ws.Range("I136").Value = 1841.875
$ws.Range("K136").Value = 5525.625
$ws.Range("M136").Value = -2975.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1623.4286
$ws.Range("I20").Value = 1675.4
$ws.Range("J20").Value = 1493.5
$ws.Range("K20").Value = 1675.4
$ws.Range("L20").Value = 1493.5
$ws.Range("M20").Value = -1428.4
$ws.Range("N20").Value = -1987.5
$ws.Range("H99").Value = 31251346
$ws.Range("I99").Value = 45455696
$ws.Range("J99").Value = 1779.8
$ws.Range("K99").Value = 45455696
$ws.Range("L99").Value = 1779.8
$ws.Range("M99").Value = -45454198
$ws.Range("N99").Value = -4775.8
$ws.Range("H107").Value = 1823.8182
$ws.Range("I107").Value = 1423.9231
$ws.Range("K107").Value = 1423.9231
$ws.Range("M107").Value = 496.0769

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1804.1818
$ws.Range("I99").Value = 1744.5
$ws.Range("K99").Value = 1744.5
$ws.Range("M99").Value = -246.5
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H126").Value = 1804.1818
$ws.Range("I126").Value = 1744.5
$ws.Range("K126").Value = 5233.5
$ws.Range("M126").Value = -2763.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 100067.8
$ws.Range("I8").Value = 100067.8
$ws.Range("K8").Value = 300203.4
$ws.Range("M8").Value = -300064.4
$ws.Range("H61").Value = 185.25
$ws.Range("J61").Value = 291.66666
$ws.Range("L61").Value = 874.9999799999999
$ws.Range("N61").Value = -1304.99998
$ws.Range("H80").Value = 4913.5
$ws.Range("I80").Value = 2894.5
$ws.Range("J80").Value = 5250
$ws.Range("K80").Value = 8683.5
$ws.Range("L80").Value = 15750
$ws.Range("M80").Value = -7747.5
$ws.Range("N80").Value = -17622
$ws.Range("H83").Value = 4913.5
$ws.Range("I83").Value = 2894.5
$ws.Range("J83").Value = 5250
$ws.Range("K83").Value = 26050.5
$ws.Range("L83").Value = 47250
$ws.Range("M83").Value = -21370.5
$ws.Range("N83").Value = -56610
$ws.Range("H104").Value = 4926.923
$ws.Range("J104").Value = 4999.8
$ws.Range("L104").Value = 14999.4
$ws.Range("N104").Value = -20241.4
$ws.Range("H122").Value = 432.0909
$ws.Range("I122").Value = 256.85715
$ws.Range("J122").Value = 738.75
$ws.Range("K122").Value = 2311.71435
$ws.Range("L122").Value = 6648.75
$ws.Range("M122").Value = 138.2856500000003
$ws.Range("N122").Value = -11548.75
$ws.Range("H129").Value = 18117854
$ws.Range("I129").Value = 83334810
$ws.Range("J129").Value = 4387969
$ws.Range("K129").Value = 250004430
$ws.Range("L129").Value = 13163907
$ws.Range("M129").Value = -249999430
$ws.Range("N129").Value = -13173907
$ws.Range("H140").Value = 2133.3215
$ws.Range("I140").Value = 1860.3334
$ws.Range("J140").Value = 2952.2856
$ws.Range("K140").Value = 5581.0002
$ws.Range("L140").Value = 8856.856800000001
$ws.Range("M140").Value = -401.0002000000004
$ws.Range("N140").Value = -19216.8568

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 18755944
$ws.Range("I70").Value = 17862372
$ws.Range("J70").Value = 20006942
$ws.Range("K70").Value = 17862372
$ws.Range("L70").Value = 20006942
$ws.Range("M70").Value = -17862102
$ws.Range("N70").Value = -20007482
$ws.Range("H73").Value = 18755944
$ws.Range("I73").Value = 17862372
$ws.Range("J73").Value = 20006942
$ws.Range("K73").Value = 17862372
$ws.Range("L73").Value = 20006942
$ws.Range("M73").Value = -17861436
$ws.Range("N73").Value = -20008814
$ws.Range("H126").Value = 2291.6667
$ws.Range("I126").Value = 1839.25
$ws.Range("K126").Value = 5517.75
$ws.Range("M126").Value = -3047.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2431.75
$ws.Range("I7").Value = 3000
$ws.Range("K7").Value = 3000
$ws.Range("M7").Value = -2888
$ws.Range("H16").Value = 1314.7142
$ws.Range("I16").Value = 1254.3334
$ws.Range("K16").Value = 1254.3334
$ws.Range("M16").Value = -1084.3334
$ws.Range("H40").Value = 4374.75
$ws.Range("I40").Value = 4999
$ws.Range("K40").Value = 4999
$ws.Range("M40").Value = -4863
$ws.Range("H46").Value = 2501
$ws.Range("I46").Value = 2001
$ws.Range("K46").Value = 2001
$ws.Range("M46").Value = -1813
$ws.Range("H61").Value = 1260.75
$ws.Range("I61").Value = 1302.7778
$ws.Range("J61").Value = 1134.6666
$ws.Range("K61").Value = 1302.7778
$ws.Range("L61").Value = 1134.6666
$ws.Range("M61").Value = -1100.7778
$ws.Range("N61").Value = -1538.6666
$ws.Range("H113").Value = 1260.75
$ws.Range("I113").Value = 1302.7778
$ws.Range("J113").Value = 1134.6666
$ws.Range("K113").Value = 1302.7778
$ws.Range("L113").Value = 1134.6666
$ws.Range("M113").Value = 867.2221999999999
$ws.Range("N113").Value = -5474.6666
$ws.Range("H126").Value = 2431.75
$ws.Range("I126").Value = 3000
$ws.Range("K126").Value = 9000
$ws.Range("M126").Value = -6530

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 57275068
$ws.Range("I122").Value = 60002404
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 180007212
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -180004762
$ws.Range("N122").Value = -7900
$ws.Range("H126").Value = 142858560
$ws.Range("I126").Value = 200000980
$ws.Range("J126").Value = 2500.5
$ws.Range("K126").Value = 600002940
$ws.Range("L126").Value = 7501.5
$ws.Range("M126").Value = -600000470
$ws.Range("N126").Value = -12441.5
